$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4564253333333334
$ws.Range("H2").Value = 1.369276
$ws.Range("I2").Value = 0.009290313076622009
$ws.Range("J2").Value = 0.009290313076622009
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.926994
$ws.Range("N2").Value = 2.780982
$ws.Range("O2").Value = 0.002566096653125693
$ws.Range("P2").Value = 0.002566096653125693
$ws.Range("Q2").Value = 0.4231035454480001
$ws.Range("R2").Value = 3.807931909032
$ws.Range("S2").Value = [double]"2.38398412924096E-05"
$ws.Range("T2").Value = [double]"2.38398412924096E-05"
$ws.Range("G3").Value = 0.4564253333333334
$ws.Range("H3").Value = 1.369276
$ws.Range("I3").Value = 0.009290313076622009
$ws.Range("J3").Value = 0.009290313076622009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 93.12610233333334
$ws.Range("N3").Value = 279.378307
$ws.Range("O3").Value = 0.2577908589658698
$ws.Range("P3").Value = 0.2577908589658698
$ws.Range("Q3").Value = 42.50511229952578
$ws.Range("R3").Value = 382.546010695732
$ws.Range("S3").Value = 0.00239495778808424
$ws.Range("T3").Value = 0.00239495778808424
$ws.Range("G4").Value = 0.4564253333333334
$ws.Range("H4").Value = 1.369276
$ws.Range("I4").Value = 0.009290313076622009
$ws.Range("J4").Value = 0.009290313076622009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 264.9957936666667
$ws.Range("N4").Value = 794.9873809999999
$ws.Range("O4").Value = 0.733559029746061
$ws.Range("P4").Value = 0.733559029746061
$ws.Range("Q4").Value = 120.9507934562396
$ws.Range("R4").Value = 1088.557141106156
$ws.Range("S4").Value = 0.006814993046523984
$ws.Range("T4").Value = 0.006814993046523984
$ws.Range("G5").Value = 0.4564253333333334
$ws.Range("H5").Value = 1.369276
$ws.Range("I5").Value = 0.009290313076622009
$ws.Range("J5").Value = 0.009290313076622009
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.197830333333333
$ws.Range("N5").Value = 6.593490999999999
$ws.Range("O5").Value = 0.006084014634943477
$ws.Range("P5").Value = 0.006084014634943477
$ws.Range("Q5").Value = 1.003145442501778
$ws.Range("R5").Value = 9.028308982516
$ws.Range("S5").Value = [double]"5.652240072137506E-05"
$ws.Range("T5").Value = [double]"5.652240072137506E-05"
$ws.Range("G6").Value = 12.034054
$ws.Range("H6").Value = 36.102162
$ws.Range("I6").Value = 0.2449472478323772
$ws.Range("J6").Value = 0.2449472478323772
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.926994
$ws.Range("N6").Value = 2.780982
$ws.Range("O6").Value = 0.002566096653125693
$ws.Range("P6").Value = 0.002566096653125693
$ws.Range("Q6").Value = 11.155495853676
$ws.Range("R6").Value = 100.399462683084
$ws.Range("S6").Value = 0.0006285583128550128
$ws.Range("T6").Value = 0.0006285583128550129
$ws.Range("G7").Value = 12.034054
$ws.Range("H7").Value = 36.102162
$ws.Range("I7").Value = 0.2449472478323772
$ws.Range("J7").Value = 0.2449472478323772
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 93.12610233333334
$ws.Range("N7").Value = 279.378307
$ws.Range("O7").Value = 0.2577908589658698
$ws.Range("P7").Value = 0.2577908589658698
$ws.Range("Q7").Value = 1120.684544288859
$ws.Range("R7").Value = 10086.16089859973
$ws.Range("S7").Value = 0.0631451614200343
$ws.Range("T7").Value = 0.06314516142003432
$ws.Range("G8").Value = 12.034054
$ws.Range("H8").Value = 36.102162
$ws.Range("I8").Value = 0.2449472478323772
$ws.Range("J8").Value = 0.2449472478323772
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 264.9957936666667
$ws.Range("N8").Value = 794.9873809999999
$ws.Range("O8").Value = 0.733559029746061
$ws.Range("P8").Value = 0.733559029746061
$ws.Range("Q8").Value = 3188.973690757524
$ws.Range("R8").Value = 28700.76321681772
$ws.Range("S8").Value = 0.1796832654588866
$ws.Range("T8").Value = 0.1796832654588866
$ws.Range("G9").Value = 12.034054
$ws.Range("H9").Value = 36.102162
$ws.Range("I9").Value = 0.2449472478323772
$ws.Range("J9").Value = 0.2449472478323772
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.197830333333333
$ws.Range("N9").Value = 6.593490999999999
$ws.Range("O9").Value = 0.006084014634943477
$ws.Range("P9").Value = 0.006084014634943477
$ws.Range("Q9").Value = 26.44880891417133
$ws.Range("R9").Value = 238.039280227542
$ws.Range("S9").Value = 0.00149026264060131
$ws.Range("T9").Value = 0.00149026264060131
$ws.Range("G10").Value = 20.326383
$ws.Range("H10").Value = 60.97914900000001
$ws.Range("I10").Value = 0.4137335243997426
$ws.Range("J10").Value = 0.4137335243997426
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.926994
$ws.Range("N10").Value = 2.780982
$ws.Range("O10").Value = 0.002566096653125693
$ws.Range("P10").Value = 0.002566096653125693
$ws.Range("Q10").Value = 18.842435082702
$ws.Range("R10").Value = 169.581915744318
$ws.Range("S10").Value = 0.001061680212248077
$ws.Range("T10").Value = 0.001061680212248077
$ws.Range("G11").Value = 20.326383
$ws.Range("H11").Value = 60.97914900000001
$ws.Range("I11").Value = 0.4137335243997426
$ws.Range("J11").Value = 0.4137335243997426
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 93.12610233333334
$ws.Range("N11").Value = 279.378307
$ws.Range("O11").Value = 0.2577908589658698
$ws.Range("P11").Value = 0.2577908589658698
$ws.Range("Q11").Value = 1892.916823324527
$ws.Range("R11").Value = 17036.25140992074
$ws.Range("S11").Value = 0.1066567206379863
$ws.Range("T11").Value = 0.1066567206379863
$ws.Range("G12").Value = 20.326383
$ws.Range("H12").Value = 60.97914900000001
$ws.Range("I12").Value = 0.4137335243997426
$ws.Range("J12").Value = 0.4137335243997426
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 264.9957936666667
$ws.Range("N12").Value = 794.9873809999999
$ws.Range("O12").Value = 0.733559029746061
$ws.Range("P12").Value = 0.733559029746061
$ws.Range("Q12").Value = 5386.405995457641
$ws.Range("R12").Value = 48477.65395911877
$ws.Range("S12").Value = 0.3034979627320935
$ws.Range("T12").Value = 0.3034979627320935
$ws.Range("G13").Value = 20.326383
$ws.Range("H13").Value = 60.97914900000001
$ws.Range("I13").Value = 0.4137335243997426
$ws.Range("J13").Value = 0.4137335243997426
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.197830333333333
$ws.Range("N13").Value = 6.593490999999999
$ws.Range("O13").Value = 0.006084014634943477
$ws.Range("P13").Value = 0.006084014634943477
$ws.Range("Q13").Value = 44.673941124351
$ws.Range("R13").Value = 402.065470119159
$ws.Range("S13").Value = 0.002517160817414778
$ws.Range("T13").Value = 0.002517160817414778
$ws.Range("G14").Value = 1.967920333333333
$ws.Range("H14").Value = 5.903761
$ws.Range("I14").Value = 0.04005605007284947
$ws.Range("J14").Value = 0.04005605007284947
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.926994
$ws.Range("N14").Value = 2.780982
$ws.Range("O14").Value = 0.002566096653125693
$ws.Range("P14").Value = 0.002566096653125693
$ws.Range("Q14").Value = 1.824250341478
$ws.Range("R14").Value = 16.418253073302
$ws.Range("S14").Value = 0.0001027876960293742
$ws.Range("T14").Value = 0.0001027876960293742
$ws.Range("G15").Value = 1.967920333333333
$ws.Range("H15").Value = 5.903761
$ws.Range("I15").Value = 0.04005605007284947
$ws.Range("J15").Value = 0.04005605007284947
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 93.12610233333334
$ws.Range("N15").Value = 279.378307
$ws.Range("O15").Value = 0.2577908589658698
$ws.Range("P15").Value = 0.2577908589658698
$ws.Range("Q15").Value = 183.2647503458475
$ws.Range("R15").Value = 1649.382753112627
$ws.Range("S15").Value = 0.01032608355505975
$ws.Range("T15").Value = 0.01032608355505975
$ws.Range("G16").Value = 1.967920333333333
$ws.Range("H16").Value = 5.903761
$ws.Range("I16").Value = 0.04005605007284947
$ws.Range("J16").Value = 0.04005605007284947
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 264.9957936666667
$ws.Range("N16").Value = 794.9873809999999
$ws.Range("O16").Value = 0.733559029746061
$ws.Range("P16").Value = 0.733559029746061
$ws.Range("Q16").Value = 521.490610604438
$ws.Range("R16").Value = 4693.415495439941
$ws.Range("S16").Value = 0.02938347722689909
$ws.Range("T16").Value = 0.02938347722689909
$ws.Range("G17").Value = 1.967920333333333
$ws.Range("H17").Value = 5.903761
$ws.Range("I17").Value = 0.04005605007284947
$ws.Range("J17").Value = 0.04005605007284947
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.197830333333333
$ws.Range("N17").Value = 6.593490999999999
$ws.Range("O17").Value = 0.006084014634943477
$ws.Range("P17").Value = 0.006084014634943477
$ws.Range("Q17").Value = 4.325155002183444
$ws.Range("R17").Value = 38.926395019651
$ws.Range("S17").Value = 0.0002437015948612449
$ws.Range("T17").Value = 0.0002437015948612449
$ws.Range("G18").Value = 2.948245666666667
$ws.Range("H18").Value = 8.844737
$ws.Range("I18").Value = 0.06001008986528831
$ws.Range("J18").Value = 0.06001008986528831
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.926994
$ws.Range("N18").Value = 2.780982
$ws.Range("O18").Value = 0.002566096653125693
$ws.Range("P18").Value = 0.002566096653125693
$ws.Range("Q18").Value = 2.733006043526
$ws.Range("R18").Value = 24.597054391734
$ws.Range("S18").Value = 0.0001539916907570884
$ws.Range("T18").Value = 0.0001539916907570884
$ws.Range("G19").Value = 2.948245666666667
$ws.Range("H19").Value = 8.844737
$ws.Range("I19").Value = 0.06001008986528831
$ws.Range("J19").Value = 0.06001008986528831
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 93.12610233333334
$ws.Range("N19").Value = 279.378307
$ws.Range("O19").Value = 0.2577908589658698
$ws.Range("P19").Value = 0.2577908589658698
$ws.Range("Q19").Value = 274.5586276578066
$ws.Range("R19").Value = 2471.027648920259
$ws.Range("S19").Value = 0.01547005261299171
$ws.Range("T19").Value = 0.01547005261299171
$ws.Range("G20").Value = 2.948245666666667
$ws.Range("H20").Value = 8.844737
$ws.Range("I20").Value = 0.06001008986528831
$ws.Range("J20").Value = 0.06001008986528831
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 264.9957936666667
$ws.Range("N20").Value = 794.9873809999999
$ws.Range("O20").Value = 0.733559029746061
$ws.Range("P20").Value = 0.733559029746061
$ws.Range("Q20").Value = 781.2727003626441
$ws.Range("R20").Value = 7031.454303263797
$ws.Range("S20").Value = 0.04402094329655482
$ws.Range("T20").Value = 0.04402094329655482
$ws.Range("G21").Value = 2.948245666666667
$ws.Range("H21").Value = 8.844737
$ws.Range("I21").Value = 0.06001008986528831
$ws.Range("J21").Value = 0.06001008986528831
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 2.197830333333333
$ws.Range("N21").Value = 6.593490999999999
$ws.Range("O21").Value = 0.006084014634943477
$ws.Range("P21").Value = 0.006084014634943477
$ws.Range("Q21").Value = 6.479743756318555
$ws.Range("R21").Value = 58.317693806867
$ws.Range("S21").Value = 0.0003651022649846873
$ws.Range("T21").Value = 0.0003651022649846873
$ws.Range("G22").Value = 11.39613766666667
$ws.Range("H22").Value = 34.188413
$ws.Range("I22").Value = 0.2319627747531205
$ws.Range("J22").Value = 0.2319627747531205
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.926994
$ws.Range("N22").Value = 2.780982
$ws.Range("O22").Value = 0.002566096653125693
$ws.Range("P22").Value = 0.002566096653125693
$ws.Range("Q22").Value = 10.564151240174
$ws.Range("R22").Value = 95.07736116156599
$ws.Range("S22").Value = 0.0005952388999437316
$ws.Range("T22").Value = 0.0005952388999437315
$ws.Range("G23").Value = 11.39613766666667
$ws.Range("H23").Value = 34.188413
$ws.Range("I23").Value = 0.2319627747531205
$ws.Range("J23").Value = 0.2319627747531205
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 93.12610233333334
$ws.Range("N23").Value = 279.378307
$ws.Range("O23").Value = 0.2577908589658698
$ws.Range("P23").Value = 0.2577908589658698
$ws.Range("Q23").Value = 1061.277882550755
$ws.Range("R23").Value = 9551.500942956791
$ws.Range("S23").Value = 0.05979788295171351
$ws.Range("T23").Value = 0.0597978829517135
$ws.Range("G24").Value = 11.39613766666667
$ws.Range("H24").Value = 34.188413
$ws.Range("I24").Value = 0.2319627747531205
$ws.Range("J24").Value = 0.2319627747531205
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 264.9957936666667
$ws.Range("N24").Value = 794.9873809999999
$ws.Range("O24").Value = 0.733559029746061
$ws.Range("P24").Value = 0.733559029746061
$ws.Range("Q24").Value = 3019.928545712928
$ws.Range("R24").Value = 27179.35691141635
$ws.Range("S24").Value = 0.1701583879851032
$ws.Range("T24").Value = 0.1701583879851032
$ws.Range("G25").Value = 11.39613766666667
$ws.Range("H25").Value = 34.188413
$ws.Range("I25").Value = 0.2319627747531205
$ws.Range("J25").Value = 0.2319627747531205
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 2.197830333333333
$ws.Range("N25").Value = 6.593490999999999
$ws.Range("O25").Value = 0.006084014634943477
$ws.Range("P25").Value = 0.006084014634943477
$ws.Range("Q25").Value = 25.04677704664255
$ws.Range("R25").Value = 225.4209934197829
$ws.Range("S25").Value = 0.001411264916360082
$ws.Range("T25").Value = 0.001411264916360082
